$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A/B values in rows 2-4 (decimal data)
$ws.Range("A2").Value = 0.8121888540000001
$ws.Range("B2").Value = 0.2739612258

$ws.Range("A3").Value = 0.2138739478333333
$ws.Range("B3").Value = 0.2538263863333334

$ws.Range("A4").Value = 0.537878838
$ws.Range("B4").Value = 0.7541248822857144

# Update column A values (integer category labels) in rows 5, 7-14, 16, 18-23
$ws.Range("A5").Value = 2

$ws.Range("A7").Value = 2
$ws.Range("A8").Value = 3
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 3
$ws.Range("A11").Value = 1
$ws.Range("A12").Value = 2
$ws.Range("A13").Value = 3
$ws.Range("A14").Value = 2

$ws.Range("A16").Value = 2

$ws.Range("A18").Value = 3
$ws.Range("A19").Value = 3
$ws.Range("A20").Value = 1
$ws.Range("A21").Value = 2
$ws.Range("A22").Value = 1
$ws.Range("A23").Value = 2
